$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: E1 label "Cost" -> "Cost per Unit" ---
$ws.Range("E1").Value = "Cost per Unit"

# --- F2:F10 "Digikey Cart" -> specific Digikey product links ---
$ws.Range("F2").Value  = "https://www.digikey.com/en/products/detail/cal-chip-electronics-inc/GMC21X7R103J50NT/22577164"
$ws.Range("F3").Value  = "https://www.digikey.com/en/products/detail/comchip-technology/CDBW46-G/3308556"
$ws.Range("F4").Value  = "https://www.digikey.com/en/products/detail/würth-elektronik/156125RS57000/9857918"
$ws.Range("F5").Value  = "https://www.digikey.com/en/products/detail/vishay-general-semiconductor-diodes-division/TZMC3V3-GS08/3104272"
$ws.Range("F6").Value  = "https://www.digikey.com/en/products/detail/keystone-electronics/7790/2171010"
$ws.Range("F7").Value  = "https://www.digikey.com/en/products/detail/bourns-inc/CRM2512-FX-51R0ELF/4698407"
$ws.Range("F8").Value  = "https://www.digikey.com/en/products/detail/yageo/RC0805FR-071KL/727444"
$ws.Range("F9").Value  = "https://www.digikey.com/en/products/detail/c-k/JS202011AQN/1640096"
$ws.Range("F10").Value = "https://www.digikey.com/en/products/detail/c-k/JS102011SAQN/1640095"

# --- Row 14 (PCB): add "Get quote at JLCPCB" link in F14 ---
$ws.Range("F14").Value = "Get quote at JLCPCB"

# --- Row 15: Red Wire -> 24AWG Red Enamel Wire, now with Qty + Amazon link, cost shown as "-" ---
$ws.Range("A15").Value = "24AWG Red Enamel Wire"
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "https://www.amazon.com/BNTECHGO-AWG-Magnet-Wire-Transformers/dp/B07DYHWLN4/ref=sr_1_5?dib=eyJ2IjoiMSJ9.T3Ul6AAFydmV-d_-CiWqSOo4gGqiz7pPGX4uu-tgp91Uz8UDmDqH8WcqGDNCB0nASmRUuhfX4iN0RZakJU07DzAEplyXYHcojTMP4m3FGJ5JPHVgh9WsCtBBpnL-l5gtozaZ1Zccg2B47I56QOWE9TdXpyhM0Vq-aou5cLVF9d2If9T7cT1eBVhsQfHPZ0qK1IjkOvqHsRwgkZCcfU0RhXIqun1H3sMrp-vOuHH7bHs.keJAopHhCjqzf_ORW6_Log1OwxKHoD8gj3qxIPjLaFU&dib_tag=se&keywords=enameled%2Bcopper%2Bwire&qid=1766271692&sr=8-5&th=1"

# --- Row 16: Green Wire -> 26AWG Green Enamel Wire, now with Qty + Amazon link, cost shown as "-" ---
$ws.Range("A16").Value = "26AWG Green Enamel Wire"
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "https://www.amazon.com/BNTECHGO-AWG-Magnet-Wire-Transformers/dp/B07HRKKMPS/ref=sr_1_7_sspa?crid=12PVFVFD45WC4&dib=eyJ2IjoiMSJ9.GfGCp1Wvx140pjLuMz3aI6wPrbZKYdBhCtRZ81xGu0XV6En0x5usAIXu_J4MJaXPsP8uD7tngx8WYsTls0IkxnFTFDnihgmKaGvu-rIKCmCyTX7a91kWyUFfrHj6u5lWP90B42q8rxUwtCxCaqo-L1I87Z-IyXnA_l0naAeAXoElYpJOL3EA33wngyAXwEBi08_0SWSdL1r4gnfuAwbsb6yLQ70njgJULhZOC5WoDEk.dYYUVAtyiKCv2QUEBzB8P9SjXI8Gwk2B3BbFux3aPSo&dib_tag=se&keywords=enameled%2Bwire%2Bgreen&qid=1766271862&sprefix=enameled%2Bwire%2Bgreen%2Caps%2C223&sr=8-7-spons&sp_csd=d2lkZ2V0TmFtZT1zcF9tdGY&th=1"

# --- New row 22: 3D Print Files ---
$ws.Range("A22").Value = "3D Print Files"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 'See "3D Printables" Folder'
$ws.Range("E22").Value = "-"
$ws.Range("F22").Value = 'See "3D Printables" Folder'
$ws.Range("F22").HorizontalAlignment = -4131

# --- Row 24 (Total): add footnote in F24 ---
$ws.Range("F24").Value = " (per board at volume, total cost may vary minimum quantity buy and shipping)"

# --- Column widths (D, E, F) ---
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668
$ws.Columns.Item(5).ColumnWidth = 14.666666666666666
$ws.Columns.Item(6).ColumnWidth = 100.83333333333334

# --- Force recalculation so SUM(E2:E21) reflects the new "-" text entries ---
$excel.Calculate()

# --- Update the saved selection/active cell to F22 ---
$ws.Range("F22").Select() | Out-Null
